# Rename the inline picture shapes that live in the document's headers
# and footers. The document has a distinct "first page" header/footer
# (w:titlePg is set), so the three logo pictures are:
#   - Footers(wdHeaderFooterPrimary)   -> footer1.xml -> PearsonLogo.png
#   - Footers(wdHeaderFooterFirstPage) -> footer2.xml -> PearsonLogo.png
#   - Headers(wdHeaderFooterFirstPage) -> header2.xml -> BTec_Logo-Orange
#
# Helper: given a Range (a header's or footer's Range), find the
# InlineShape it contains and rename it. Some header/footer ranges have
# several paragraphs before the one holding the picture, so we look
# paragraph-by-paragraph rather than indexing the shape straight off the
# whole range.
function Rename-InlineShapeInRange($rng, $newName) {
    $paras = $rng.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.InlineShapes.Count -gt 0) {
            $shp = $p.Range.InlineShapes(1)
            $shp.Name = $newName
        }
    }
}

$d = $word.ActiveDocument
$section = $d.Sections(1)

# footer1.xml (default/primary footer) : image1.png -> image2.png
Rename-InlineShapeInRange $section.Footers(1).Range "image2.png"

# footer2.xml (first-page footer) : image1.png -> image2.png
Rename-InlineShapeInRange $section.Footers(2).Range "image2.png"

# header2.xml (first-page header) : image2.jpg -> image1.jpg
Rename-InlineShapeInRange $section.Headers(2).Range "image1.jpg"
